$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.583.71"
$ws.Range("E2").Value = "  +1.83%  "

$ws.Range("D3").Value = "2.615.53"
$ws.Range("E3").Value = "  +1.23%  "

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.79%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.46%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  +1.93%  "

$ws.Range("D9").Value = "2.612.85"
$ws.Range("E9").Value = "  +1.20%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.127"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +11.86%  "

$ws.Range("E11").Value = "  +0.79%  "

$ws.Range("E12").Value = "  +1.17%  "

$ws.Range("E13").Value = "  +0.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.13%  "

$ws.Range("E15").Value = "  +4.43%  "

$ws.Range("D16").Value = "3.092.06"
$ws.Range("E16").Value = "  +1.51%  "

$ws.Range("D17").Value = "67.596.04"
$ws.Range("E17").Value = "  +2.09%  "

$ws.Range("D18").Value = "2.617.64"
$ws.Range("E18").Value = "  +1.36%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.26%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "362.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.16%  "

$ws.Range("E21").Value = "  -1.84%  "

$ws.Range("E22").Value = "  -0.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.07%  "

$ws.Range("E24").Value = "  -0.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.82%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.30%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000106"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.19%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "582.75"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.54%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.17%  "

$ws.Range("E31").Value = "  -0.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.95"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.31%  "

$ws.Range("E33").Value = "  +1.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.130"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.88%  "

$ws.Range("E35").Value = "  -0.03%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.53"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.51%  "

$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.42"
$ws.Range("D38").Style = "Normal"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "156.06"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.09%  "

$ws.Range("E40").Value = "  +1.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.40"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.15%  "

$ws.Range("E42").Value = "  +3.89%  "

$ws.Range("E43").Value = "  +4.09%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.27%  "

$ws.Range("E45").Value = "  +0.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.01%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "156.87"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.87%  "

$ws.Range("D48").Value = "0.0₆0291"
$ws.Range("E48").Value = "  -5.51%  "

$ws.Range("E49").Value = "  +0.80%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.19%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.622"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.48%  "
